# The commit swaps the two embedded themes: the presentation's main
# theme (ppt/theme/theme1.xml, used by the slide master / all slides)
# changes from the custom "Integral" theme to the built-in
# "Office Theme" palette (this is what the notes master's theme,
# theme2.xml, already contained). The dark1/light1 colors (black/white)
# and the font scheme are identical between the two themes, so only
# the ten theme colors that differ need to be updated.
#
# PowerPoint RGB color values are packed as 0x00BBGGRR (red in the low
# byte), i.e. value = R + G*256 + B*65536.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# MsoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
